$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 4 new problem rows (180-183) to the tracking sheet, mirroring the
# existing row layout/format (column styles: A/D/E/F/G/H/I = style 1 (or 4
# for dates), B/C = style 2). We copy an existing row's formatting first so
# the cells pick up the already-defined styles instead of minting new ones.
# ---------------------------------------------------------------------------

function Add-Row($RowNum, $Height, $A, $B, $C, $D, $E, $F, $G, $H, $I) {
    $srcRow = $RowNum - 1
    $ws.Range("A" + $srcRow + ":I" + $srcRow).Copy()
    $ws.Range("A" + $RowNum).PasteSpecial(-4122)

    $ws.Range("A" + $RowNum).Value = $A
    $ws.Range("B" + $RowNum).Value = $B
    $ws.Range("C" + $RowNum).Value = $C
    $ws.Range("D" + $RowNum).Value = $D
    $ws.Range("E" + $RowNum).Value = $E
    $ws.Range("F" + $RowNum).Value = $F
    if ($null -ne $G) {
        $ws.Range("G" + $RowNum).Value = $G
    }
    $ws.Range("H" + $RowNum).Value = $H
    $ws.Range("I" + $RowNum).Value = $I

    $ws.Rows.Item($RowNum).RowHeight = $Height
}

Add-Row 180 51 3494 "Find the Minimum Amount of Time to Brew Potions" "#dynamic-programming " "medium" 0 1 35 45939 45939
Add-Row 181 17 8 "String to Integer (atoi)" "#math #string" "medium" 0 1 20 45943 45943
Add-Row 182 34 3349 "Adjacent Increasing Subarrays Detection I" "#array" "easy" 0 1 $null 45944 45944
Add-Row 183 34 3350 "Adjacent Increasing Subarrays Detection II" "#array" "medium" 1 0 5 45945 45945

# Match the final selection/scroll state recorded in the workbook.
$ws.Range("C182").Select() | Out-Null

Write-Host "done"
